$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.173.76"
$ws.Range("E2").Value = "  +0.69%  "
# Row 3
$ws.Range("D3").Value = "1.848.03"
$ws.Range("E3").Value = "  +1.11%  "
# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.40%  "
# Row 5
$ws.Range("D5").Value = "'312.86"
$ws.Range("E5").Value = "  +0.32%  "
# Row 6
$ws.Range("E6").Value = "  -0.30%  "
# Row 7
$ws.Range("D7").Value = "'0.4604"
$ws.Range("E7").Value = "  -0.23%  "
# Row 8
$ws.Range("E8").Value = "  +0.14%  "
# Row 9
$ws.Range("D9").Value = "'0.07270"
$ws.Range("E9").Value = "  -0.90%  "
# Row 10
$ws.Range("D10").Value = "'0.8828"
$ws.Range("E10").Value = "  +0.98%  "
# Row 11
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.980.65"
$ws.Range("E11").Value = "  +7.49%  "
# Row 12
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'19.97"
$ws.Range("E12").Value = "  +0.79%  "
# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07810"
$ws.Range("E13").Value = "  -1.61%  "
# Row 14
$ws.Range("D14").Value = "'5.365"
$ws.Range("E14").Value = "  +0.46%  "
# Row 15
$ws.Range("D15").Value = "'6.494"
$ws.Range("E15").Value = "  -0.80%  "
# Row 16
$ws.Range("D16").Value = "'91.22"
$ws.Range("E16").Value = "  -0.17%  "
# Row 17
$ws.Range("E17").Value = "  -0.38%  "
# Row 18
$ws.Range("D18").Value = "'0.000008906"
$ws.Range("E18").Value = "  +0.23%  "
# Row 20
$ws.Range("E20").Value = "  -0.64%  "
# Row 21
$ws.Range("D21").Value = "27.209.38"
$ws.Range("E21").Value = "  +0.15%  "
# Row 22
$ws.Range("D22").Value = "'5.043"
$ws.Range("E22").Value = "  -1.36%  "
# Row 23
$ws.Range("D23").Value = "'10.48"
$ws.Range("E23").Value = "  -0.61%  "
# Row 24
$ws.Range("D24").Value = "2.153.14"
$ws.Range("E24").Value = "  +3.03%  "
# Row 25
$ws.Range("D25").Value = "'1.956"
$ws.Range("E25").Value = "  +5.78%  "
# Row 26
$ws.Range("D26").Value = "'151.73"
$ws.Range("E26").Value = "  -0.80%  "
# Row 27
$ws.Range("E27").Value = "  -0.04%  "
# Row 28
$ws.Range("E28").Value = "  +1.15%  "
# Row 29
$ws.Range("D29").Value = "'115.35"
$ws.Range("E29").Value = "  -0.06%  "
# Row 30
$ws.Range("D30").Value = "'5.053"
$ws.Range("E30").Value = "  -1.80%  "
# Row 31
$ws.Range("D31").Value = "'0.08823"
$ws.Range("E31").Value = "  -0.90%  "
# Row 32
$ws.Range("D32").Value = "'3.095"
$ws.Range("E32").Value = "  +4.53%  "
# Row 33
$ws.Range("D33").Value = "'0.7603"
# Row 34
$ws.Range("E34").Value = "  +3.42%  "
# Row 35
$ws.Range("D35").Value = "'4.491"
$ws.Range("E35").Value = "  +1.42%  "
# Row 36
$ws.Range("D36").Value = "'2.719"
$ws.Range("E36").Value = "  +9.90%  "
# Row 37
$ws.Range("D37").Value = "'1.080"
$ws.Range("E37").Value = "  +0.90%  "
# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05231"
$ws.Range("E38").Value = "  +0.09%  "
# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01938"
$ws.Range("E39").Value = "  -0.84%  "
# Row 40
$ws.Range("D40").Value = "'2.937"
$ws.Range("E40").Value = "  -0.19%  "
# Row 41
$ws.Range("D41").Value = "'7.061"
$ws.Range("E41").Value = "  -0.69%  "
# Row 42
$ws.Range("D42").Value = "'0.5088"
$ws.Range("E42").Value = "  -1.47%  "
# Row 43
$ws.Range("D43").Value = "'0.1625"
$ws.Range("E43").Value = "  -0.09%  "
# Row 44
$ws.Range("D44").Value = "'8.362"
$ws.Range("E44").Value = "  +2.11%  "
# Row 45
$ws.Range("E45").Value = "  -1.41%  "
# Row 46
$ws.Range("D46").Value = "'10.35"
$ws.Range("E46").Value = "  +1.81%  "
# Row 47
$ws.Range("E47").Value = "  -0.34%  "
# Row 48
$ws.Range("D48").Value = "'102.56"
$ws.Range("E48").Value = "  +0.01%  "
# Row 49
$ws.Range("D49").Value = "'1.632"
$ws.Range("E49").Value = "  -0.05%  "
# Row 50
$ws.Range("E50").Value = "  +0.29%  "
# Row 51
$ws.Range("D51").Value = "'65.56"
$ws.Range("E51").Value = "  +1.27%  "
